# Modificar vista de pagina de carga de plantillas de paginación.
#
# The template sheet "Planificacion" gets a horario code ("HR01") filled
# in for the week's schedule columns (B5, and C5:C12), and the active
# selection moves from H15 to E15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planificacion")

# Fill the new HR01 schedule values.
$ws.Range("B5").Value = "HR01"
$ws.Range("C5").Value = "HR01"
$ws.Range("C6").Value = "HR01"
$ws.Range("C7").Value = "HR01"
$ws.Range("C8").Value = "HR01"
$ws.Range("C9").Value = "HR01"
$ws.Range("C10").Value = "HR01"
$ws.Range("C11").Value = "HR01"
$ws.Range("C12").Value = "HR01"

# Move the active selection/cell to E15 (was H15).
$ws.Range("E15").Select()
